# Regen save_data: recalculated "K" column (col G) values replacing the
# previous Strike# derived numbers, after regenerating std/mean and the
# s_vals used to compute them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(0,5,2,3,2,0,5,6,4,4,3,5,3,1,2,5,4,3,3,3,5,2,1,2,1,1,1,1,1,1,0,0,0,1,0,1,0,1,0,1,4,3,2,2)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
